$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - set values, then copy H1's format (bold/border/alignment) onto I1:J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for columns I (I0) and J (IF)
$data = @(
    @(4, 4),
    @(4, 6),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(4, 4),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(1, 3),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(6, 8),
    @(5, 6),
    @(1, 2),
    @(1, 2),
    @(1, 2),
    @(1, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
